$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.566.07'
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").Value = '3.466.90'
$ws.Range("E3").Value = '  +2.31%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''577.33'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").Value = '''162.10'
$ws.Range("E6").Value = '  +4.77%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '3.469.88'
$ws.Range("E8").Value = '  +2.31%  '
$ws.Range("E9").Value = '  +10.16%  '
$ws.Range("D10").Value = '''7.36'
$ws.Range("E10").Value = '  -1.68%  '
$ws.Range("E11").Value = '  +5.22%  '
$ws.Range("E12").Value = '  +2.47%  '
$ws.Range("D13").Value = '4.064.05'
$ws.Range("E13").Value = '  +2.32%  '
$ws.Range("E14").Value = '  -2.71%  '
$ws.Range("D15").Value = '''0.0000196'
$ws.Range("E15").Value = '  +6.09%  '
$ws.Range("D16").Value = '''29.05'
$ws.Range("E16").Value = '  +7.54%  '
$ws.Range("D17").Value = '64.581.32'
$ws.Range("E17").Value = '  +1.83%  '
$ws.Range("D18").Value = '3.473.30'
$ws.Range("E18").Value = '  +2.52%  '
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("D20").Value = '''14.58'
$ws.Range("E20").Value = '  +4.38%  '
$ws.Range("D21").Value = '''392.63'
$ws.Range("E21").Value = '  +0.91%  '
$ws.Range("D22").Value = '''8.27'
$ws.Range("E22").Value = '  -1.76%  '
$ws.Range("E23").Value = '  +3.23%  '
$ws.Range("D24").Value = '''73.42'
$ws.Range("E24").Value = '  +3.71%  '
$ws.Range("E25").Value = '  +0.48%  '
$ws.Range("D26").Value = '''0.0000125'
$ws.Range("E26").Value = '  +21.04%  '
$ws.Range("D27").Value = '''9.57'
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D30").Value = '''6.20'
$ws.Range("E30").Value = '  +11.05%  '
$ws.Range("E31").Value = '  +9.89%  '
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("D33").Value = '''6.59'
$ws.Range("E33").Value = '  +2.25%  '
$ws.Range("E34").Value = '  +2.94%  '
$ws.Range("D35").Value = '''0.999'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").Value = '''7.14'
$ws.Range("E36").Value = '  +6.11%  '
$ws.Range("E37").Value = '  +2.37%  '
$ws.Range("D38").Value = '''161.34'
$ws.Range("E38").Value = '  +1.67%  '
$ws.Range("E39").Value = '  +0.66%  '
$ws.Range("E40").Value = '  +3.67%  '
$ws.Range("D41").Value = '''27.68'
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").Value = '2.925.37'
$ws.Range("E42").Value = '  +1.70%  '
$ws.Range("E43").Value = '  +6.97%  '
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("D45").Value = '''42.83'
$ws.Range("E45").Value = '  +4.06%  '
$ws.Range("D46").Value = '''0.777'
$ws.Range("E46").Value = '  +1.99%  '
$ws.Range("D47").Value = '''24.22'
$ws.Range("E47").Value = '  +9.68%  '
$ws.Range("E48").Value = '  +3.79%  '
$ws.Range("D49").Value = '''2.22'
$ws.Range("E49").Value = '  +15.43%  '
$ws.Range("D50").Value = '''0.878'
$ws.Range("E50").Value = '  +8.15%  '
$ws.Range("E51").Value = '  +4.75%  '
